$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tbl8")

# Drop the unused quantile columns (I:P) and the per-predictor detail rows (11:20);
# the report now only keeps the 0.1/0.2/0.5 quantiles plus model-fit summary rows.
$ws.Columns("I:P").Delete()
$ws.Rows("11:20").Delete()

# Clear cells that must end up blank in the new 8x10 table
# (one call per cell -- this host only honours the first area of a comma union range)
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("B9").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("F9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("B10").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("F10").ClearContents()
$ws.Range("H10").ClearContents()

# Write the refreshed coefficient/AIC/BIC/chi-square table (values kept as text,
# matching the source report export, via a leading apostrophe + style reset)
$ws.Range("A2").Value = "'1"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "'ss1"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "'0.15 (-0.39, 0.71)"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "'0.6"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.06 (0.3, 1.85)"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "'0"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "'1.28 (0.33, 2.19)"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = "'0"
$ws.Range("H2").Style = "Normal"
$ws.Range("A3").Value = "'2"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "'ss2"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "'-0.23 (-0.76, 0.27)"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "'0.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-0.2 (-0.8, 0.36)"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = "'0.52"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "'-0.04 (-0.73, 0.68)"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value = "'0.92"
$ws.Range("H3").Style = "Normal"
$ws.Range("A4").Value = "'3"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = "'ss3"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "'-0.53 (-1.02, 0.05)"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'0.04"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.12 (-0.66, 0.43)"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "'0.68"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = "'-0.22 (-0.8, 0.31)"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = "'0.45"
$ws.Range("H4").Style = "Normal"
$ws.Range("A5").Value = "'4"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = "'ss4"
$ws.Range("B5").Style = "Normal"
$ws.Range("E5").Value = "'-0.24 (-0.86, 0.31)"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "'0.41"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = "'-0.28 (-0.84, 0.28)"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = "'0.37"
$ws.Range("H5").Style = "Normal"
$ws.Range("A6").Value = "'5"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "'ss5"
$ws.Range("B6").Style = "Normal"
$ws.Range("E6").Value = "'-0.84 (-1.58, -0.03)"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = "'0.03"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = "'-0.1 (-0.76, 0.54)"
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value = "'0.77"
$ws.Range("H6").Style = "Normal"
$ws.Range("A7").Value = "'6"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = "'ss6"
$ws.Range("B7").Style = "Normal"
$ws.Range("G7").Value = "'-0.93 (-1.74, -0.04)"
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = "'0.03"
$ws.Range("H7").Style = "Normal"
$ws.Range("A8").Value = "'aic"
$ws.Range("A8").Style = "Normal"
$ws.Range("C8").Value = "'1312.563"
$ws.Range("C8").Style = "Normal"
$ws.Range("E8").Value = "'1305.808"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'1308.856"
$ws.Range("G8").Style = "Normal"
$ws.Range("A9").Value = "'bic"
$ws.Range("A9").Style = "Normal"
$ws.Range("C9").Value = "'1373.397"
$ws.Range("C9").Style = "Normal"
$ws.Range("E9").Value = "'1374.246"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'1381.096"
$ws.Range("G9").Style = "Normal"
$ws.Range("A10").Value = "'pr_chisq"
$ws.Range("A10").Style = "Normal"
$ws.Range("E10").Value = "'0.005"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'1"
$ws.Range("G10").Style = "Normal"
